# Add a second worksheet "majeurs+mineurs" after "Sheet" with a lookup
# table, point a VLOOKUP at it from the main sheet, and give H13 an
# explicit MAX formula to match its already-cached value.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- create the new sheet, positioned after "Sheet" -----------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "majeurs+mineurs"

# --- header row -------------------------------------------------------------
$ws2.Range("A1").Value = "Noma"
$ws2.Range("B1").Value = "Nom"
$ws2.Range("C1").Value = "Etude"
$ws2.Range("D1").Value = "Résultat"
$ws2.Range("E1").Value = "Note finale"

# --- row 2 : Smith, Adam -----------------------------------------------------
$ws2.Range("A2").NumberFormat = "@"
$ws2.Range("A2").Value = "00000001"
$ws1.Range("E12").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$ws2.Range("B2").Value = "Smith, Adam"

$ws2.Range("C2").Value = "OSIS2MA"
$ws1.Range("D12").Copy()
$ws2.Range("C2").PasteSpecial(-4122)

$ws2.Range("D2").Value = 14.7
$ws2.Range("E2").Formula = "=ROUND(D2,0)"

# --- row 3 : Doe, John --------------------------------------------------------
$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = "00000002"
$ws1.Range("E13").Copy()
$ws2.Range("A3").PasteSpecial(-4122)

$ws2.Range("B3").Value = "Doe, John"

$ws2.Range("C3").Value = "OSIS2MA"
$ws1.Range("D13").Copy()
$ws2.Range("C3").PasteSpecial(-4122)

$ws2.Range("D3").Value = 17.1
$ws2.Range("E3").Formula = "=ROUND(D3,0)"

# --- point the main sheet's formulas at the new table ------------------------
$ws1.Range("H12").Formula = "=VLOOKUP(E12,'majeurs+mineurs'!`$A`$2:`$E`$3,5,0)"
$ws1.Range("H13").Formula = "=MAX(1,12,17)"

# --- restore the active selection on the main sheet ---------------------------
$ws1.Activate()
$ws1.Range("E34").Select() | Out-Null
